$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("epoch50", "epoch100", "epoch150", "epoch200", "epoch250", "epoch300", "epoch350", "epoch400")
$values  = @(16.79890526307596, 13.649135908565, 12.14348361701578, 11.22498699539417, 10.99800219809687, 10.80918658424068, 10.70239408998876, 10.6191528004569)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2  # starting at column B
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
